$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# Mark the 2016 / 2022 / 2023 commemorative varieties as owned (0 -> 1)
$ws.Range("F12").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1

# Move the cursor/selection on the frozen bottom-right pane to J14
$ws.Activate()
$ws.Range("J14").Select()
